# Generate Report for Handback
# Adds a new handback record (6c2306a1-8b8d-488a-9da4-cb910a9f0a4f) as row 3
# to the "Overview", "zh-cn" and "de-de" worksheets/tables.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Overview sheet (table3 / "Overview")
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md"
$wsOverview.Range("B3").Value = "e2e\6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-11-14 06:07:24"
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b9ed5aed37893fe994f2e6a7258144151a3743e8/e2e/6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md",
    [Type]::Missing,
    [Type]::Missing,
    "e2e\6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md"
) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# zh-cn sheet (table1 / "zh_cn")
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3").Value = "6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = "6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.3539ced2b77af5979eadfe2a3a0e91a6fc3b2f4a.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-11-14 06:07:12"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("I3").Value = "6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md"
$wsZhCn.Range("J3").Value = "6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.3539ced2b77af5979eadfe2a3a0e91a6fc3b2f4a.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-11-14 06:08:06"
$wsZhCn.Range("K3").NumberFormat = $dateFmt
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f7279f5f7138a732713a8c33ea24fedf71a758e7/e2e/6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md",
    [Type]::Missing,
    [Type]::Missing,
    "6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f7279f5f7138a732713a8c33ea24fedf71a758e7/e2e/6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md",
    [Type]::Missing,
    [Type]::Missing,
    "6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md"
) | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------------
# de-de sheet (table2 / "de_de")
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3").Value = "6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = "6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.3539ced2b77af5979eadfe2a3a0e91a6fc3b2f4a.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-11-14 06:07:24"
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("I3").Value = "6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md"
$wsDeDe.Range("J3").Value = "6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.3539ced2b77af5979eadfe2a3a0e91a6fc3b2f4a.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-11-14 06:08:24"
$wsDeDe.Range("K3").NumberFormat = $dateFmt
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/95fcd1ac05c9025899e5764fb620baa9907a59d4/e2e/6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md",
    [Type]::Missing,
    [Type]::Missing,
    "6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/95fcd1ac05c9025899e5764fb620baa9907a59d4/e2e/6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md",
    [Type]::Missing,
    [Type]::Missing,
    "6c2306a1-8b8d-488a-9da4-cb910a9f0a4f.md"
) | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

"Handback row added to Overview, zh-cn and de-de sheets"
